# Fix model sheet -- 'type' now looks to match the prompt type and expands
# recursively.
#
# The "model" sheet's `type` column used to hold generic storage kinds
# (string/number/object) plus a separate `elementType` column (column C)
# that carried the real prompt type for compound fields (e.g. geopoint).
# Now `type` directly mirrors the prompt type used on the "survey" sheet
# (text/decimal/geopoint/select_one/integer), and the now-redundant
# `elementType` column is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# name                      | type (old -> new)
# refrigerator_id           | string   -> text
# refrigerator_size         | number   -> decimal
# refrigerator_location     | object   -> geopoint (was in elementType/col C)
# refrigerator_condition    | string   -> select_one
# refrigerator_stock_level  | integer  -> integer (unchanged)
$ws.Range("B2").Value = "text"
$ws.Range("B3").Value = "decimal"
$ws.Range("B4").Value = "geopoint"
$ws.Range("B5").Value = "select_one"

# Column C ("elementType") is no longer needed -- its only values already
# moved into column B above.
$ws.Range("C1").ClearContents()
$ws.Range("C4").ClearContents()

# Match the author's recorded selection on the model sheet.
$ws.Activate()
$ws.Range("B4").Select()
